$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("THURSDAY SINGLES")
$ws.Range("A1").Value = "test"
